# Apply a cyclic swap of the Id (A), Antal (I), Ost (Q) and Nord (R)
# values among rows 12, 14, 15, 16, 17, 19 and 20 on the active sheet,
# matching the upstream "automatic update" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, taken from the target diff.
# Column I ("Antal") holds numeric-looking text, so it is written with a
# leading apostrophe to keep it stored as text rather than being
# auto-converted to a number.

$ws.Range("A12").Value = 111378856
$ws.Range("I12").Value = "'10"
$ws.Range("Q12").Value = 505494.3524330241
$ws.Range("R12").Value = 6913043.848162009

$ws.Range("A14").Value = 111378933
$ws.Range("I14").Value = "'25"
$ws.Range("Q14").Value = 505597.6535686332
$ws.Range("R14").Value = 6913018.009825628

$ws.Range("A15").Value = 111378866
$ws.Range("I15").Value = "'10"
$ws.Range("Q15").Value = 505492.5216403615
$ws.Range("R15").Value = 6913025.731493607

$ws.Range("A16").Value = 111378893
$ws.Range("I16").Value = "'25"
$ws.Range("Q16").Value = 505612.5119866763
$ws.Range("R16").Value = 6913033.361683531

$ws.Range("A17").Value = 111378874
$ws.Range("I17").Value = "'50"
$ws.Range("Q17").Value = 505592.4968292552
$ws.Range("R17").Value = 6913042.152801346

$ws.Range("A19").Value = 111378946
$ws.Range("I19").Value = "'100"
$ws.Range("Q19").Value = 505602.791734456
$ws.Range("R19").Value = 6913005.013642685

$ws.Range("A20").Value = 111378954
$ws.Range("I20").Value = "'15"
$ws.Range("Q20").Value = 505590.6913760683
$ws.Range("R20").Value = 6913009.17353364
